$wb = $excel.ActiveWorkbook

function Set-RowValues($SheetName, $Row, $Values) {
    # $Values order: H, I, J, K, L, M, N ; $null means the cell should be cleared/absent
    $ws = $wb.Worksheets.Item($SheetName)
    $cols = @("H","I","J","K","L","M","N")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$Row"
        $ws.Range($addr).Value = $Values[$i]
    }
}

Set-RowValues "ALC" 64 @(5592.5, 4777.5, 6000, 4777.5, 6000, -4529.5, -6496)
Set-RowValues "ALC" 67 @(5592.5, 4777.5, 6000, 4777.5, 6000, -3919.5, -7716)
Set-RowValues "ALC" 76 @(83429944, 108434.2, 500037500, 108434.2, 500037500, -108119.2, -500038130)
Set-RowValues "ALC" 79 @(83429944, 108434.2, 500037500, 108434.2, 500037500, -107342.2, -500039684)
Set-RowValues "ALC" 92 @(2351.375, 902.5, 2834.3333, 902.5, 2834.3333, 345.5, -5330.3333)
Set-RowValues "ALC" 103 @(492.57144, 213.85715, 771.2857, 641.5714499999999, 2313.8571, -55.57144999999991, -3485.8571)
Set-RowValues "ALC" 112 @(2688.394, 979, 2741.8125, 2937, 8225.4375, -1829, -10441.4375)
Set-RowValues "ALC" 116 @(4192.7144, 3299.6667, 4862.5, 3299.6667, 4862.5, 142.3332999999998, -11746.5)
Set-RowValues "ALC" 127 @(11950.833, 14734.777, 3599, 44204.331, 10797, -39244.331, -20717)
Set-RowValues "ALC" 132 @(2014.5151, 1020, 6489.8335, 3060, 19469.5005, -530, -24529.5005)
Set-RowValues "ALC" 137 @(4614.41, 2801.6538, 8239.923000000001, 8404.9614, 24719.769, -5854.9614, -29819.769)
Set-RowValues "ARM" 23 @(3670666.8, 3670666.8, 0, 3670666.8, 0, -3670407.8, $null)
Set-RowValues "ARM" 32 @(2541.1357, 2212.4736, 7536.8, 2212.4736, 7536.8, -1925.4736, -8110.8)
Set-RowValues "ARM" 45 @(7147.625, 3818.5557, 11427.857, 3818.5557, 11427.857, -3441.5557, -12181.857)
Set-RowValues "ARM" 102 @(2574.75, 2574.75, 0, 2574.75, 0, -952.75, $null)
Set-RowValues "BSM" 134 @(2874.95, 1657.75, 7743.75, 4973.25, 23231.25, -2438.25, -28301.25)
Set-RowValues "CRP" 14 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "CRP" 99 @(5202.3125, 4026.3333, 6714.2856, 4026.3333, 6714.2856, -2528.3333, -9710.285599999999)
Set-RowValues "CRP" 126 @(5202.3125, 4026.3333, 6714.2856, 12078.9999, 20142.8568, -9608.999899999999, -25082.8568)
Set-RowValues "CRP" 132 @(5847.3, 4882.6665, 6636.5454, 14647.9995, 19909.6362, -12117.9995, -24969.6362)
Set-RowValues "CRP" 141 @(70332.664, 0, 70332.664, 0, 70332.664, $null, -80692.664)
Set-RowValues "CUL" 12 @(227.53334, 14, 242.78572, 42, 728.35716, 131, -1074.35716)
Set-RowValues "CUL" 38 @(33.22222, 30.666666, 38.333332, 91.99999800000001, 114.999996, 255.000002, -808.999996)
Set-RowValues "CUL" 80 @(69933.336, 0, 69933.336, 0, 209800.008, $null, -211672.008)
Set-RowValues "CUL" 82 @(5989.4443, 4046.6667, 6960.8335, 12140.0001, 20882.5005, -11734.0001, -21694.5005)
Set-RowValues "CUL" 83 @(69933.336, 0, 69933.336, 0, 629400.024, $null, -638760.024)
Set-RowValues "CUL" 85 @(5989.4443, 4046.6667, 6960.8335, 12140.0001, 20882.5005, -10736.0001, -23690.5005)
Set-RowValues "CUL" 86 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "CUL" 89 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "CUL" 92 @(2001620.4, 10000002, 2025, 30000006, 6075, -29998758, -8571)
Set-RowValues "CUL" 97 @(758.6667, 1303, 649.8, 3909, 1949.4, -3413, -2941.4)
Set-RowValues "CUL" 98 @(1599.1666, 2803.3333, 1197.7778, 8409.999899999999, 3593.3334, -6911.999899999999, -6589.3334)
Set-RowValues "GSM" 80 @(1004548.9, 913986.8, 1115235.9, 913986.8, 1115235.9, -912988.8, -1117231.9)
Set-RowValues "GSM" 83 @(1004548.9, 913986.8, 1115235.9, 4569934, 5576179.5, -4564942, -5586163.5)
Set-RowValues "GSM" 122 @(6879.7856, 6415.737, 7859.4443, 19247.211, 23578.3329, -16797.211, -28478.3329)
Set-RowValues "LTW" 4 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "LTW" 28 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "LTW" 37 @(0, 0, 0, 0, 0, $null, $null)
Set-RowValues "WVR" 96 @(113263.22, 127046.125, 3000, 127046.125, 3000, -125673.125, -5746)
Set-RowValues "WVR" 127 @(46991.25, 0, 46991.25, 0, 46991.25, $null, -56911.25)

Write-Output "Applied all row updates"